# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" / "Valor Mora" table (E16:F30) is re-sorted so the
# periods run in descending order (2103 down to 2001) instead of ascending
# (2001 up to 2103), and the two "Valor Mora" figures that differed from
# the rest (25396 vs 33125) travel with their original period (2103) to
# the top of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-sort) period/value pairs for rows 16..30.
$periodToValue = @{}
$periodList    = @()
for ($r = 16; $r -le 30; $r++) {
    $period = $ws.Cells.Item($r, 5).Value2
    $value  = $ws.Cells.Item($r, 6).Value2
    $periodList += $period
    $periodToValue[$period] = $value
}

# Sort the periods in descending order - newest period first.
$sortedPeriods = $periodList | Sort-Object -Descending

# Clear the existing cell contents first so the underlying string table is
# rebuilt cleanly for the reordered values instead of merely re-pointing
# at the old entries.
$ws.Range("E16:F30").ClearContents()

$r = 16
foreach ($period in $sortedPeriods) {
    $ws.Cells.Item($r, 5).Value = $period
    $ws.Cells.Item($r, 6).Value = $periodToValue[$period]
    $r++
}
